$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Plot")
$ws2 = $wb.Worksheets.Item("Droplet diameters")

# --- Droplet diameters sheet: re-run of the detector after adding a blur
# step to cope with noisy backgrounds -> slightly different droplet sizes,
# and the extra spurious tiny droplets (rows 9-12) are no longer detected.
$ws2.Range("A2").Value = 23.55065224750685
$ws2.Range("B2").Value = 6.030617155121099
$ws2.Range("A3").Value = 19.9056400039448
$ws2.Range("A4").Value = 9.996371919901266
$ws2.Range("A5").Value = 2.4392042757347
$ws2.Range("B5").Value = 2.088255084412402
$ws2.Range("A6").Value = 1.32030739689701
$ws2.Range("A7").Value = 0.9231165248209556
$ws2.Range("A8").Value = 0.07456865054232387

# Rows 9-12 no longer exist -> clearing them shrinks the sheet back to A1:B8
$ws2.Range("A9:A12").ClearContents() | Out-Null

# --- Plot sheet: frequency/average/max formulas now look at A2:A8 instead
# of the old A2:A12 range.
$ws1.Range("C2:C5").FormulaArray = "=frequency('Droplet diameters'!A2:A8,'Plot'!A2:A4)"
$ws1.Range("C7").Formula = "=AVERAGE('Droplet diameters'!A2:A8)"
$ws1.Range("C8").Formula = "=MAX('Droplet diameters'!A2:A8)"
